$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    6.119791666666667,
    6.119791666666667,
    6.119791666666667,
    6.119791666666667,
    6.119791666666667,
    6.25,
    6.25,
    6.25,
    6.25,
    6.380208333333334,
    6.380208333333334,
    6.380208333333334,
    6.510416666666667,
    6.510416666666667,
    6.510416666666667,
    6.510416666666667,
    6.640625000000001,
    6.640625000000001,
    6.640625000000001,
    6.640625000000001,
    6.640625000000001,
    6.640625000000001,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.770833333333334,
    6.901041666666667,
    6.901041666666667,
    6.901041666666667,
    7.031250000000001,
    7.031250000000001,
    7.031250000000001,
    7.161458333333334,
    7.161458333333334,
    7.161458333333334,
    7.161458333333334,
    7.161458333333334,
    7.161458333333334,
    7.161458333333334,
    7.291666666666667,
    7.291666666666667,
    7.291666666666667,
    7.291666666666667,
    7.291666666666667,
    7.291666666666667,
    7.291666666666667,
    7.421875000000001,
    7.552083333333334,
    8.072916666666668,
    8.984375,
    9.765625,
    9.895833333333334,
    9.895833333333334,
    9.895833333333334,
    9.375,
    10.02604166666667,
    10.02604166666667,
    10.02604166666667,
    10.02604166666667,
    10.15625,
    10.41666666666667,
    10.80729166666667,
    10.28645833333333,
    9.635416666666668,
    9.635416666666668,
    12.5,
    12.5,
    12.5,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333,
    12.63020833333333
)

$startRow = 151
$endRow = $startRow + $values.Length - 1

$range = $ws.Range("A$startRow`:A$endRow")
$data = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $data[$i,0] = $values[$i]
}
$range.Value = $data
